$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item("SPRINT 0")

# Duplicate "SPRINT 0" right after itself -> becomes the new "SPRINT 1" sheet,
# inheriting layout/formatting (column width, date style, page margins, ...).
$ws0.Copy($null, $ws0)
$ws1 = $wb.Worksheets.Item($ws0.Index + 1)
$ws1.Name = "SPRINT 1"

# --- SPRINT 1 content ---
# Row 2 (Date / Quoi / Temps (h)) is already correct from the copy.

# Row 3
$ws1.Range("B3").Value = 43904
$ws1.Range("C3").Value = "Lecture TOJC"
$ws1.Range("D3").Value = 1

# Row 4
$ws1.Range("B4").Value = 43904
$ws1.Range("C4").Value = "Lecture rapport analyse 2020 WP1"
$ws1.Range("D4").Value = 2

# Row 5
$ws1.Range("B5").Value = 43911
$ws1.Range("C5").Value = "Synthèse état de l'art - introduction"
$ws1.Range("D5").Value = 0.5

# Row 6
$ws1.Range("B6").Value = 43911
$ws1.Range("C6").Value = "Synthèse état de l'art - GNI + Trust Project"
$ws1.Range("D6").Value = 1

# Row 7
$ws1.Range("B7").Value = 43911
$ws1.Range("C7").Value = "Synthèse état de l'art - JTI + Transparency Journalism + NewsGuard"
$ws1.Range("D7").Value = 2

# Row 8
$ws1.Range("B8").Value = 43911
$ws1.Range("C8").Value = "Synthèse état de l'art - comparaison"
$ws1.Range("D8").Value = 2

# Row 9
$ws1.Range("B9").Value = 43911
$ws1.Range("C9").Value = "Synthèse état de l'art - conclusion"
$ws1.Range("D9").Value = 2

# Rows 10-14 used to hold extra SPRINT 0 entries; SPRINT 1 only keeps the
# (empty, date-formatted) B cells.
$ws1.Range("C10:D14").ClearContents()
$ws1.Range("B10").Value = ""
$ws1.Range("B11").Value = ""
$ws1.Range("B12").Value = ""
$ws1.Range("B13").Value = ""
$ws1.Range("B14").Value = ""

# Row 15 held the totals formula in SPRINT 0; in SPRINT 1 the total moves to
# row 16 and only sums rows 3:9.
$ws1.Range("B15:D15").ClearContents()
$ws1.Range("D16").Formula = "=SUM(D3:D9)"

# Row 17 footer note
$ws1.Range("C17").Value = "Sprint 1 review, notes de séance en document annexe"

# --- Selection / active tab bookkeeping ---
# Select on SPRINT 0 first (leaves it not tab-selected), then select on the
# new SPRINT 1 last so it ends up the active sheet/tab.
$ws0.Range("C18").Select()
$ws1.Range("D10").Select()
